$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the data (all columns except the row-index column A) between
#    row 168 <-> row 169, and between row 173 <-> row 174.
# ---------------------------------------------------------------------------
$pairs = @(168,169),@(173,174)
foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$($r1):AC$($r1)")
    $range2 = $ws.Range("B$($r2):AC$($r2)")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# ---------------------------------------------------------------------------
# 2) Append a new row (188) with the latest match for this league.
# ---------------------------------------------------------------------------

# Column A carries the same bold/centered/bordered style used by the rest of
# column A - copy formats only from the row above, then set the value.
$ws.Range("A187").Copy() | Out-Null
$ws.Range("A188").PasteSpecial(-4122) | Out-Null
$ws.Range("A188").Value = 186

$ws.Range("B188").Value = 7793507
$ws.Range("C188").Value = "Chile Primera B"
$ws.Range("D188").Value = "Chile Primera B"

# Column E carries the custom date/time number format - copy formats only
# from the row above, then set the value.
$ws.Range("E187").Copy() | Out-Null
$ws.Range("E188").PasteSpecial(-4122) | Out-Null
$ws.Range("E188").Value = 45387.79166666666

$ws.Range("F188").Value = "Deportes Limache"
$ws.Range("G188").Value = "Curico Unido"

# H188, I188 (FTHG/FTAG) and J188 (FTR) are left blank - the match hasn't
# been played yet, so those fields are not populated.

$ws.Range("K188").Value = 1.8
$ws.Range("L188").Value = 3.4
$ws.Range("M188").Value = 4
$ws.Range("N188").Value = 1.95
$ws.Range("O188").Value = 3.3
$ws.Range("P188").Value = 3.5
$ws.Range("Q188").Value = -0.5
$ws.Range("R188").Value = 2
$ws.Range("S188").Value = 1.8
$ws.Range("T188").Value = 2.5
$ws.Range("U188").Value = 1.95
$ws.Range("V188").Value = 1.85
$ws.Range("W188").Value = 0
$ws.Range("X188").Value = 0
$ws.Range("Y188").Value = 0
$ws.Range("Z188").Value = 0
$ws.Range("AA188").Value = 0

# AB188, AC188 (PL_AhOver/PL_AhUnder) are also left blank for the same reason.
